$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 3 and 4 (Lenah Cheloti <-> Ochieng Charles)
$ws.Range("A3").Value = "Ochieng Charles"
$ws.Range("B3").Value = "'1.00"
$ws.Range("C3").Value = "'30.00"
$ws.Range("D3").Value = "'-29.00"
$ws.Range("E3").Value = "'3.33%"

$ws.Range("A4").Value = "Lenah Cheloti"
$ws.Range("B4").Value = "'3.00"
$ws.Range("C4").Value = "'22.00"
$ws.Range("D4").Value = "'-19.00"
$ws.Range("E4").Value = "'13.64%"

# Insert a new row before the totals row (current row 5) to hold Moses Ngugi
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "Moses  Ngugi"
$ws.Range("B5").Value = "'0.00"
$ws.Range("C5").Value = "'25.00"
$ws.Range("D5").Value = "'-25.00"
$ws.Range("E5").Value = "'0.00%"

# Update totals row (now row 6)
$ws.Range("A6").Value = "KD Totals"
$ws.Range("B6").Value = "'7.00"
$ws.Range("C6").Value = "'102.00"
$ws.Range("D6").Value = "'-95.00"
$ws.Range("E6").Value = "'7.24%"
